# "Generate Report for Handback"
#
# The localization status for the two source files (bc176e68... and
# be632df2...) flips from "In Translation" to "Handed back: in sync with
# en-US", the Latest Target File / Latest Handback File columns on the
# per-language sheets get populated (with hyperlinks, mirroring the
# existing "Source File Name" hyperlink look) and the Latest Handback
# DateTime placeholder gets a real timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"

$zh.Range("I2").Value = "bc176e68-2f81-4b06-9ed9-06b73a8d066b.md"
$zh.Range("J2").Value = "bc176e68-2f81-4b06-9ed9-06b73a8d066b.98b41539b7bae0aefc78daa8a44b467369b3de2b.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-27 16:25:12"

$zh.Range("I3").Value = "be632df2-8ce2-40e9-b56c-63a600e8d8ba.md"
$zh.Range("J3").Value = "be632df2-8ce2-40e9-b56c-63a600e8d8ba.b256fa2589a81413d1f5bf8c5a459486c4985159.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-27 16:25:12"

$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/267916c60b29a5cbe4d07e92de27ce4737d3dfe9/e2e/bc176e68-2f81-4b06-9ed9-06b73a8d066b.md", "", "", "bc176e68-2f81-4b06-9ed9-06b73a8d066b.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/267916c60b29a5cbe4d07e92de27ce4737d3dfe9/e2e/be632df2-8ce2-40e9-b56c-63a600e8d8ba.md", "", "", "be632df2-8ce2-40e9-b56c-63a600e8d8ba.md")

$zh.Columns.Item(3).ColumnWidth = 29.144371396019366
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"

$de.Range("I2").Value = "bc176e68-2f81-4b06-9ed9-06b73a8d066b.md"
$de.Range("J2").Value = "bc176e68-2f81-4b06-9ed9-06b73a8d066b.98b41539b7bae0aefc78daa8a44b467369b3de2b.de-de.xlf"
$de.Range("K2").Value = "2016-08-27 16:25:18"

$de.Range("I3").Value = "be632df2-8ce2-40e9-b56c-63a600e8d8ba.md"
$de.Range("J3").Value = "be632df2-8ce2-40e9-b56c-63a600e8d8ba.b256fa2589a81413d1f5bf8c5a459486c4985159.de-de.xlf"
$de.Range("K3").Value = "2016-08-27 16:25:18"

$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/267916c60b29a5cbe4d07e92de27ce4737d3dfe9/e2e/bc176e68-2f81-4b06-9ed9-06b73a8d066b.md", "", "", "bc176e68-2f81-4b06-9ed9-06b73a8d066b.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/267916c60b29a5cbe4d07e92de27ce4737d3dfe9/e2e/be632df2-8ce2-40e9-b56c-63a600e8d8ba.md", "", "", "be632df2-8ce2-40e9-b56c-63a600e8d8ba.md")

$de.Columns.Item(3).ColumnWidth = 29.144371396019366
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated"
